$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 388; this shifts rows 388:424 down to 389:425
# and pushes the used range down to row 425 (dimension A1:R425).
$ws.Rows("388:388").Insert()

# Populate the newly inserted row 388 with the new record.
# Columns A, B, C, E, F, G, I, O, R are inherited/unchanged from the row
# that used to be at 388 (now at 389), so fill them in explicitly too
# since the freshly inserted row starts out blank.
$ws.Cells.Item(388, 1).Value = 6
$ws.Cells.Item(388, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(388, 3).Value = "Metropolitana"
$ws.Cells.Item(388, 4).Value = 44769
$ws.Cells.Item(388, 5).Value = 13
$ws.Cells.Item(388, 6).Value = 100112032
$ws.Cells.Item(388, 7).Value = "Zapallo italiano"
$ws.Cells.Item(388, 8).Value = "Sin especificar"
$ws.Cells.Item(388, 9).Value = "Primera"
$ws.Cells.Item(388, 10).Value = 400
$ws.Cells.Item(388, 11).Value = 14000
$ws.Cells.Item(388, 12).Value = 15000
$ws.Cells.Item(388, 13).Value = 14425
$ws.Cells.Item(388, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(388, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(388, 16).Value = 288
$ws.Cells.Item(388, 17).Value = 50
$ws.Cells.Item(388, 18).Value = "Hortaliza"
